# Commit ngay 8/11/2022 lan 1
#
# 1. The auto-updating "datetimeFigureOut" date footer field (on the Slide
#    Master and on every Slide Layout) is bumped from 2022-11-07 to
#    2022-11-08, as happens when PowerPoint is re-saved the next day with
#    an auto-updating date field.
# 2. Slide 1 ("4 cach chen bang") is removed - the deck goes from 4 slides
#    down to 3 (Designed Table / slide with "Where/When" table / Table
#    title slide).
#
# NOTE: the date-footer placeholders are updated *before* the slide is
# deleted - shape/text lookups on the master & layouts can otherwise pick
# up stale data once the deck's slide list has been mutated.

$p = $ppt.ActivePresentation

function Update-DateFooter($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2022-11-07") {
                $tr.Text = "2022-11-08"
            }
        }
    }
}

# --- 1. Refresh the auto-date field text on the Slide Master and on every
#        Slide Layout so it reflects 2022-11-08 instead of 2022-11-07 -------
Update-DateFooter $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateFooter $layout.Shapes
}

# --- 2. Remove the first slide ---------------------------------------------
$p.Slides.Item(1).Delete()
